$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1872.9697
$ws.Range("I40").Value = 2533
$ws.Range("J40").Value = 1386.6316
$ws.Range("K40").Value = 2533
$ws.Range("L40").Value = 1386.6316
$ws.Range("M40").Value = -2358
$ws.Range("N40").Value = -1736.6316

$ws.Range("H55").Value = 1500
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 1500
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 1500
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -1928

$ws.Range("H138").Value = 11015.194
$ws.Range("I138").Value = 2334.4614
$ws.Range("J138").Value = 15921.695
$ws.Range("K138").Value = 7003.3842
$ws.Range("L138").Value = 47765.085
$ws.Range("M138").Value = -1863.3842
$ws.Range("N138").Value = -58045.085

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1538
$ws.Range("I2").Value = 1538
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1538
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -1425
$ws.Range("N2").ClearContents()

$ws.Range("H74").Value = 7410.636
$ws.Range("I74").Value = 2723.6155
$ws.Range("J74").Value = 14180.777
$ws.Range("K74").Value = 2723.6155
$ws.Range("L74").Value = 14180.777
$ws.Range("M74").Value = -1849.6155
$ws.Range("N74").Value = -15928.777

$ws.Range("H77").Value = 7410.636
$ws.Range("I77").Value = 2723.6155
$ws.Range("J77").Value = 14180.777
$ws.Range("K77").Value = 13618.0775
$ws.Range("L77").Value = 70903.88499999999
$ws.Range("M77").Value = -9250.077499999999
$ws.Range("N77").Value = -79639.88499999999

$ws.Range("H116").Value = 1538
$ws.Range("I116").Value = 1538
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1538
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 756
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1538
$ws.Range("I3").Value = 1538
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1538
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -1424
$ws.Range("N3").ClearContents()

$ws.Range("H134").Value = 48177.547
$ws.Range("I134").Value = 2493.25
$ws.Range("K134").Value = 7479.75
$ws.Range("M134").Value = -4944.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5687.943
$ws.Range("I31").Value = 5532.4585
$ws.Range("J31").Value = 6027.1816
$ws.Range("K31").Value = 5532.4585
$ws.Range("L31").Value = 6027.1816
$ws.Range("M31").Value = -5237.4585
$ws.Range("N31").Value = -6617.1816

$ws.Range("H34").Value = 5687.943
$ws.Range("I34").Value = 5532.4585
$ws.Range("J34").Value = 6027.1816
$ws.Range("K34").Value = 5532.4585
$ws.Range("L34").Value = 6027.1816
$ws.Range("M34").Value = -5330.4585
$ws.Range("N34").Value = -6431.1816

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6254.4116
$ws.Range("I5").Value = 376.92307
$ws.Range("K5").Value = 1130.76921
$ws.Range("M5").Value = -1018.76921

$ws.Range("H122").Value = 886.375
$ws.Range("I122").Value = 495.66666
$ws.Range("K122").Value = 4460.99994
$ws.Range("M122").Value = -2010.99994

$ws.Range("H132").Value = 1244.4375
$ws.Range("I132").Value = 927
$ws.Range("K132").Value = 8343
$ws.Range("M132").Value = -5813

$ws.Range("H135").Value = 6254.4116
$ws.Range("I135").Value = 376.92307
$ws.Range("K135").Value = 3392.30763
$ws.Range("M135").Value = -857.3076299999998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 143356.42
$ws.Range("I3").Value = 200299
$ws.Range("J3").Value = 1000
$ws.Range("K3").Value = 200299
$ws.Range("L3").Value = 1000
$ws.Range("M3").Value = -200183
$ws.Range("N3").Value = -1232

$ws.Range("H10").Value = 250002750
$ws.Range("I10").Value = 1000000000
$ws.Range("J10").Value = 3668
$ws.Range("K10").Value = 1000000000
$ws.Range("L10").Value = 3668
$ws.Range("M10").Value = -999999831
$ws.Range("N10").Value = -4006

$ws.Range("H11").Value = 343333340
$ws.Range("I11").Value = 343333340
$ws.Range("K11").Value = 343333340
$ws.Range("M11").Value = -343333201

$ws.Range("H120").Value = 39900
$ws.Range("J120").Value = 39900
$ws.Range("L120").Value = 39900
$ws.Range("N120").Value = -49576

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H10").Value = 7900
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 7900
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 7900
$ws.Range("M10").ClearContents()
$ws.Range("N10").Value = -8180

$ws.Range("H12").Value = 650
$ws.Range("I12").Value = 650
$ws.Range("K12").Value = 650
$ws.Range("M12").Value = -480

$ws.Range("H46").Value = 3000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 3000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 3000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -3376

$ws.Range("H61").Value = 23220.777
$ws.Range("I61").Value = 29383.143
$ws.Range("J61").Value = 1652.5
$ws.Range("K61").Value = 29383.143
$ws.Range("L61").Value = 1652.5
$ws.Range("M61").Value = -29181.143
$ws.Range("N61").Value = -2056.5

$ws.Range("H113").Value = 23220.777
$ws.Range("I113").Value = 29383.143
$ws.Range("J113").Value = 1652.5
$ws.Range("K113").Value = 29383.143
$ws.Range("L113").Value = 1652.5
$ws.Range("M113").Value = -27213.143
$ws.Range("N113").Value = -5992.5

$ws.Range("H136").Value = 6214.276
$ws.Range("I136").Value = 3591.4
$ws.Range("J136").Value = 7594.737
$ws.Range("K136").Value = 10774.2
$ws.Range("L136").Value = 22784.211
$ws.Range("M136").Value = -8224.200000000001
$ws.Range("N136").Value = -27884.211

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1868.25
$ws.Range("I13").Value = 1868.25
$ws.Range("K13").Value = 1868.25
$ws.Range("M13").Value = -1728.25

$ws.Range("H81").Value = 14288999
$ws.Range("I81").Value = 1749
$ws.Range("J81").Value = 16670208
$ws.Range("K81").Value = 3498
$ws.Range("L81").Value = 33340416
$ws.Range("M81").Value = -2437
$ws.Range("N81").Value = -33342538

$ws.Range("H84").Value = 14288999
$ws.Range("I84").Value = 1749
$ws.Range("J84").Value = 16670208
$ws.Range("K84").Value = 17490
$ws.Range("L84").Value = 166702080
$ws.Range("M84").Value = -12186
$ws.Range("N84").Value = -166712688

$ws.Range("H126").Value = 1289.5834
$ws.Range("I126").Value = 1267.6471
$ws.Range("K126").Value = 3802.9413
$ws.Range("M126").Value = -1332.9413
